$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1 with same style as the existing header (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I2:J9
$data = @(
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(3, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
